# MasterBarangMultiSatuanPerPelanggan.xlsx - "update master barang multi satuan per pelanggan"
#
# The HARGABELI1 / HARGABELI2 / HARGABELI3 columns (G, M, S) are collapsed into a
# single HARGABELI column that sits right after JMLSTOK (new column E), and the
# first product row is rewritten to an AQUA / MINUMAN entry with fresh pricing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the now-redundant HARGABELI2 / HARGABELI3 columns outright (delete from the
# right first so earlier column letters stay valid).
$ws.Columns("S").Delete()
$ws.Columns("M").Delete()

# Relocate HARGABELI1 (col G) so it becomes a single merged HARGABELI column right
# after JMLSTOK (col E); ISI1/SATUAN1 slide right to F/G.
$ws.Columns("G").Cut()
$ws.Columns("E").Insert()

# The merged column is no longer specific to "satuan 1" - rename its header.
$ws.Range("E1").Value = "HARGABELI"

# Row 2 becomes a new product: AQUA / MINUMAN, with updated quantities & pricing.
$ws.Range("A2").Value = "AQUA"
$ws.Range("B2").Value = "MINUMAN"
$ws.Range("D2").Value = 10000
$ws.Range("E2").Value = 1000
$ws.Range("H2").Value = 1500
$ws.Range("I2").Value = 1400
$ws.Range("J2").Value = 1300
$ws.Range("K2").Value = 48
$ws.Range("M2").Value = 55000
$ws.Range("N2").Value = 54000
$ws.Range("O2").Value = 53000

# Resize columns to fit the reshuffled content and restore the saved selection.
$ws.Columns.AutoFit()
$ws.Range("P2").Select()
